# Refresh Price / Volume(1h) columns with the latest feed snapshot.
# Values are stored as plain text in the sheet (e.g. "299.10", "-1.69%"),
# so each cell is forced to Text format before the write and returned to
# the sheet's default ("Normal") style afterwards to avoid leaving any
# stray number-format/style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2'  = '299.10'
    'E2'  = '-1.69%'
    'E3'  = '-1.53%'
    'D4'  = '5.087'
    'E4'  = '-2.36%'
    'D5'  = '0.07913'
    'E5'  = '1.21%'
    'D6'  = '2.321'
    'E6'  = '-2.62%'
    'D7'  = '7.811'
    'E7'  = '-2.15%'
    'D8'  = '3.870'
    'E8'  = '0.21%'
    'D9'  = '0.9234'
    'E9'  = '1.25%'
    'E10' = '0.60%'
    'D11' = '0.07584'
    'E11' = '3.13%'
    'D12' = '0.09297'
    'E12' = '14.28%'
    'D13' = '0.03010'
    'E13' = '-1.27%'
    'D14' = '0.1002'
    'E14' = '0.80%'
    'D15' = '0.001522'
    'E15' = '-0.31%'
    'D16' = '0.005908'
    'E16' = '-4.60%'
    'D17' = '3.476'
    'E17' = '-0.70%'
    'D18' = '2.266'
    'E18' = '1.09%'
    'E19' = '-0.81%'
    'E20' = '-0.15%'
    'D21' = '4.015'
    'E21' = '-13.75%'
    'D22' = '0.1703'
    'E22' = '8.96%'
    'D23' = '0.04625'
    'E23' = '-0.49%'
    'D24' = '0.001252'
    'E24' = '-0.64%'
    'D25' = '0.004475'
    'E25' = '-1.37%'
    'D26' = '0.0001251'
    'E26' = '-7.13%'
    'D27' = '0.0003401'
    'E27' = '24.28%'
    'E39' = '-3.07%'
    'D40' = '0.04598'
    'E40' = '0.07%'
    'D41' = '0.006978'
    'E41' = '-4.04%'
    'D42' = '0.1358'
    'E42' = '-0.32%'
    'D43' = '0.002193'
    'E43' = '-1.93%'
    'E44' = '-6.02%'
    'D45' = '0.00006294'
    'E45' = '-2.96%'
    'D46' = '0.00000000751'
    'E46' = '0.30%'
    'D47' = '0.007989'
    'E47' = '-19.14%'
    'D48' = '1.155'
    'E48' = '40.79%'
    'D49' = '0.00002102'
    'E49' = '0.30%'
    'D50' = '0.0002002'
    'E50' = '0.30%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}
